$d = $word.ActiveDocument

# --- "Programa" section: Portuguese bullet list ---
$ptSrc = "- Tratamento de imagens: resolução, definição, contraste, saturação; uso de técnicas automatizadas de determinação de tamanho e distribuição de partículas.- Proposição e ajuste de equações empíricas a resultados de medidas experimentais: as diversas propostas de relações para a deformação plástica e encruamento.- Potenciais interatômicos e o método de dinâmica molecular clássica; simulação de solidificação de um metal puro.- Cinética de nucleação e crescimento: a equação de Johnson-Mehl-Avrami-Kolmogorov (JMAK) e sua aplicação computacional.- Elementos finitos: estudo do estado de tensão de materiais sob carregamentos mecânicos; simulação de transferência de calor em tratamentos térmicos.- Método de Monte Carlo aplicado à transição ferro-paramagnética e à cinética de crescimento de grão- Cálculo de diagramas de fases: curvas de energia livre, o método CALPHAD; Thermo-Calc e Dictra."
$ptDst = "- Tratamento de imagens: resolução, definição, contraste, saturação; uso de técnicas automatizadas de determinação de tamanho e distribuição de partículas.^l- Proposição e ajuste de equações empíricas a resultados de medidas experimentais: as diversas propostas de relações para a deformação plástica e encruamento.^l- Potenciais interatômicos e o método de dinâmica molecular clássica; simulação de solidificação de um metal puro.^l- Cinética de nucleação e crescimento: a equação de Johnson-Mehl-Avrami-Kolmogorov (JMAK) e sua aplicação computacional.^l- Elementos finitos: estudo do estado de tensão de materiais sob carregamentos mecânicos; simulação de transferência de calor em tratamentos térmicos.^l- Método de Monte Carlo aplicado à transição ferro-paramagnética e à cinética de crescimento de grão^l- Cálculo de diagramas de fases: curvas de energia livre, o método CALPHAD; Thermo-Calc e Dictra."
$r1 = $d.Content
$found1 = $r1.Find.Execute($ptSrc, $true, $false, $false, $false, $false, $true, 1, $false, $ptDst, 2)

# --- "Programa" section: English (italic) bullet list ---
$enSrc = "- Image treatment: resolution, definition, contrast, saturation; use of automated techniques for determining particle size and distribution.- Proposition and fit of empirical equations to results of experimental measures: the various proposals for relationships for plastic deformation and hardening.- Interatomic potentials and the classical molecular dynamics method; simulation of solidification of a pure metal.- Nucleation and growth kinetics: the Johnson-Mehl-Avrami-Kolmogorov (JMAK) equation and its computational application.- Finite element method: study of the stress state of materials under mechanical loads; simulation of heat transfer applied to heat treatments.- Monte Carlo method applied to the ferro-paramagnetic transition and to grain growth kinetics- Calculation of phase diagrams: free energy curves, the CALPHAD method; Thermo-Calc and Dictra."
$enDst = "- Image treatment: resolution, definition, contrast, saturation; use of automated techniques for determining particle size and distribution.^l- Proposition and fit of empirical equations to results of experimental measures: the various proposals for relationships for plastic deformation and hardening.^l- Interatomic potentials and the classical molecular dynamics method; simulation of solidification of a pure metal.^l- Nucleation and growth kinetics: the Johnson-Mehl-Avrami-Kolmogorov (JMAK) equation and its computational application.^l- Finite element method: study of the stress state of materials under mechanical loads; simulation of heat transfer applied to heat treatments.^l- Monte Carlo method applied to the ferro-paramagnetic transition and to grain growth kinetics^l- Calculation of phase diagrams: free energy curves, the CALPHAD method; Thermo-Calc and Dictra."
$r2 = $d.Content
$found2 = $r2.Find.Execute($enSrc, $true, $false, $false, $false, $false, $true, 1, $false, $enDst, 2)

# --- "Bibliografia" section ---
$bibSrc = "- Richard LESAR, Computational Materials Science – Fundamentals to Applications. MRS, 2013.- Rob Phillips, Crystals, Defects and Microstructures – Modelling across scales. Cambridge, 2001.- Artigos publicados em revistas como Computational Materials Science, Calphad, Journal of Alloys and Compounds, etc."
$bibDst = "- Richard LESAR, Computational Materials Science – Fundamentals to Applications. MRS, 2013.^l- Rob Phillips, Crystals, Defects and Microstructures – Modelling across scales. Cambridge, 2001.^l- Artigos publicados em revistas como Computational Materials Science, Calphad, Journal of Alloys and Compounds, etc."
$r3 = $d.Content
$found3 = $r3.Find.Execute($bibSrc, $true, $false, $false, $false, $false, $true, 1, $false, $bibDst, 2)

Write-Host "Programa (PT) replaced: $found1"
Write-Host "Programa (EN) replaced: $found2"
Write-Host "Bibliografia replaced: $found3"
